# Insert a new "neighbourhood_group" column between Bydel_Navn (B) and
# Inntekt 2017 (old C, becomes D), populate it, tidy up the stray
# right-hand formatting column, and restore column widths / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 0. Clear the stray formatted-but-empty cells in column H (H3:H19) ---
# These carried no value, only left-over cell formatting, and are not part
# of the sheet any more once the new column is inserted.
$ws.Range("H3:H19").Clear() | Out-Null

# --- 1. Insert a new column before column C -----------------------------
# Shifts the existing "Inntekt 2017" column (C) one place right, to D.
$ws.Columns("C:C").Insert()

# --- 2. Header for the new column ---------------------------------------
$ws.Range("C1").Value = "neighbourhood_group"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").Font.Name = "Calibri"
$ws.Range("C1").Font.Size = 11
$ws.Range("C1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("C1").VerticalAlignment = -4160     # xlTop

# --- 3. Data values for the new column -----------------------------------
$values = @{
    "Vestre Aker" = 16
    "Gamle Oslo" = 3
    "Frogner" = 2
    "Østensjø" = 9
    "St.Hanshaugen" = 13
    "Nordre Aker" = 7
    "Sagene" = 10
    "Grünerløkka" = 5
    "Sentrum" = 11
    "Alna" = 0
    "Grorud" = 4
    "Stovner" = 14
    "Bjerke" = 1
    "Ullern" = 15
    "Nordstrand" = 8
    "Søndre Nordstrand" = 12
}

for ($r = 2; $r -le 17; $r++) {
    $name = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 3).Value = $values[$name]
}

# --- 4. Column widths -----------------------------------------------------
$ws.Columns("C:C").ColumnWidth = 22.296875
$ws.Columns("D:D").ColumnWidth = 14.09765625
$ws.Columns("F:F").ColumnWidth = 24.59765625

# --- 5. Selection, matching the final saved state -------------------------
$ws.Range("F10").Select()
